# Auto-generated script to fix number/text formatting in shared strings
# (replicates commit: "fix: fixed formatting when scrapping floating point numbers")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original cell styles for the ranges we are about to touch,
# since assigning a leading-quote (text-forcing) value can otherwise
# introduce a quote-prefix style bit on the cell.
$hRange = $ws.Range("H2:H293")
$hStyle = $hRange.Style

$s_E205 = $ws.Range("E205").Style
$s_E207 = $ws.Range("E207").Style
$s_E217 = $ws.Range("E217").Style

# --- Column E: normalize punctuation in two "Razon social" entries ---
$ws.Range("E205").Value = "'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E207").Value = "'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E217").Value = "'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- Column H: convert "Importe" values from es-AR (1.234,56) to plain (1234.56) ---
$ws.Range("H2").Value = "'27500.00"
$ws.Range("H3").Value = "'10560.00"
$ws.Range("H4").Value = "'10400.00"
$ws.Range("H5").Value = "'18200.00"
$ws.Range("H6").Value = "'6100.00"
$ws.Range("H7").Value = "'42425.00"
$ws.Range("H8").Value = "'10800.00"
$ws.Range("H9").Value = "'578500.00"
$ws.Range("H10").Value = "'2616000.00"
$ws.Range("H11").Value = "'1308000.00"
$ws.Range("H12").Value = "'1308000.00"
$ws.Range("H13").Value = "'850000.00"
$ws.Range("H14").Value = "'250000.00"
$ws.Range("H15").Value = "'5437200.00"
$ws.Range("H16").Value = "'2943000.00"
$ws.Range("H17").Value = "'155.00"
$ws.Range("H18").Value = "'3000.00"
$ws.Range("H19").Value = "'583773.94"
$ws.Range("H20").Value = "'14100.00"
$ws.Range("H21").Value = "'48400.00"
$ws.Range("H22").Value = "'479500.00"
$ws.Range("H23").Value = "'12899.90"
$ws.Range("H24").Value = "'1537.92"
$ws.Range("H25").Value = "'808956.00"
$ws.Range("H26").Value = "'701551.47"
$ws.Range("H27").Value = "'293580.32"
$ws.Range("H28").Value = "'98266.84"
$ws.Range("H29").Value = "'13641.50"
$ws.Range("H30").Value = "'138459.60"
$ws.Range("H31").Value = "'50125.00"
$ws.Range("H32").Value = "'85800.00"
$ws.Range("H33").Value = "'2400.00"
$ws.Range("H34").Value = "'54947.20"
$ws.Range("H35").Value = "'3480.46"
$ws.Range("H36").Value = "'79458.17"
$ws.Range("H37").Value = "'21686.00"
$ws.Range("H38").Value = "'42480.00"
$ws.Range("H39").Value = "'58503.89"
$ws.Range("H40").Value = "'14000.00"
$ws.Range("H41").Value = "'22700.00"
$ws.Range("H42").Value = "'13600.08"
$ws.Range("H43").Value = "'144.88"
$ws.Range("H44").Value = "'176.00"
$ws.Range("H45").Value = "'83.48"
$ws.Range("H46").Value = "'1869.63"
$ws.Range("H47").Value = "'6732145.46"
$ws.Range("H48").Value = "'74026.84"
$ws.Range("H49").Value = "'197615.00"
$ws.Range("H50").Value = "'4400.00"
$ws.Range("H51").Value = "'3311.00"
$ws.Range("H52").Value = "'76084.00"
$ws.Range("H53").Value = "'74215.00"
$ws.Range("H54").Value = "'1769.95"
$ws.Range("H55").Value = "'201991.00"
$ws.Range("H56").Value = "'37000.00"
$ws.Range("H57").Value = "'1576.61"
$ws.Range("H58").Value = "'14900.00"
$ws.Range("H59").Value = "'22999.00"
$ws.Range("H60").Value = "'2685.00"
$ws.Range("H61").Value = "'7000.00"
$ws.Range("H62").Value = "'56740.00"
$ws.Range("H63").Value = "'6900.00"
$ws.Range("H64").Value = "'8848.19"
$ws.Range("H65").Value = "'357290.23"
$ws.Range("H66").Value = "'69325.00"
$ws.Range("H67").Value = "'74.00"
$ws.Range("H68").Value = "'2658.00"
$ws.Range("H69").Value = "'154.00"
$ws.Range("H70").Value = "'7815.58"
$ws.Range("H71").Value = "'1928.50"
$ws.Range("H72").Value = "'63720.00"
$ws.Range("H73").Value = "'21939.23"
$ws.Range("H74").Value = "'940.00"
$ws.Range("H75").Value = "'3605.60"
$ws.Range("H76").Value = "'6369.07"
$ws.Range("H77").Value = "'511.21"
$ws.Range("H78").Value = "'459521.83"
$ws.Range("H79").Value = "'4200.00"
$ws.Range("H80").Value = "'1172.60"
$ws.Range("H81").Value = "'160.00"
$ws.Range("H82").Value = "'5051.08"
$ws.Range("H83").Value = "'11700.00"
$ws.Range("H84").Value = "'3450.00"
$ws.Range("H85").Value = "'21583.78"
$ws.Range("H86").Value = "'700.00"
$ws.Range("H87").Value = "'16200.00"
$ws.Range("H88").Value = "'34000.00"
$ws.Range("H89").Value = "'43800.00"
$ws.Range("H90").Value = "'12000.00"
$ws.Range("H91").Value = "'2000.00"
$ws.Range("H92").Value = "'650.00"
$ws.Range("H93").Value = "'18000.00"
$ws.Range("H94").Value = "'75000.00"
$ws.Range("H95").Value = "'10480.00"
$ws.Range("H96").Value = "'8282.00"
$ws.Range("H97").Value = "'7480.00"
$ws.Range("H98").Value = "'36000.00"
$ws.Range("H99").Value = "'726.00"
$ws.Range("H100").Value = "'29121.00"
$ws.Range("H101").Value = "'2174.00"
$ws.Range("H102").Value = "'9.00"
$ws.Range("H103").Value = "'6304.00"
$ws.Range("H104").Value = "'43.05"
$ws.Range("H105").Value = "'55485.36"
$ws.Range("H106").Value = "'21151.10"
$ws.Range("H107").Value = "'52955.66"
$ws.Range("H108").Value = "'7822.54"
$ws.Range("H109").Value = "'23779.00"
$ws.Range("H110").Value = "'504.00"
$ws.Range("H111").Value = "'4576.00"
$ws.Range("H112").Value = "'94938.60"
$ws.Range("H113").Value = "'360.00"
$ws.Range("H114").Value = "'90.00"
$ws.Range("H115").Value = "'420.00"
$ws.Range("H116").Value = "'1280.00"
$ws.Range("H117").Value = "'6800.00"
$ws.Range("H118").Value = "'35020.00"
$ws.Range("H119").Value = "'38081.00"
$ws.Range("H120").Value = "'21010.00"
$ws.Range("H121").Value = "'9960.00"
$ws.Range("H122").Value = "'17561.10"
$ws.Range("H123").Value = "'109588.74"
$ws.Range("H124").Value = "'2186.00"
$ws.Range("H125").Value = "'38960.00"
$ws.Range("H126").Value = "'1165.00"
$ws.Range("H127").Value = "'1449.90"
$ws.Range("H128").Value = "'40255.00"
$ws.Range("H129").Value = "'10558.00"
$ws.Range("H130").Value = "'351000.00"
$ws.Range("H131").Value = "'9426.00"
$ws.Range("H132").Value = "'67800.00"
$ws.Range("H133").Value = "'6790.00"
$ws.Range("H134").Value = "'41154.00"
$ws.Range("H135").Value = "'16990.00"
$ws.Range("H136").Value = "'13828.00"
$ws.Range("H137").Value = "'3051.20"
$ws.Range("H138").Value = "'905.20"
$ws.Range("H139").Value = "'173.40"
$ws.Range("H140").Value = "'26500.00"
$ws.Range("H141").Value = "'12544.00"
$ws.Range("H142").Value = "'7000.00"
$ws.Range("H143").Value = "'450.00"
$ws.Range("H144").Value = "'5000.00"
$ws.Range("H145").Value = "'105000.00"
$ws.Range("H146").Value = "'5000.00"
$ws.Range("H147").Value = "'145000.00"
$ws.Range("H148").Value = "'88000.00"
$ws.Range("H149").Value = "'29800.00"
$ws.Range("H150").Value = "'5000.00"
$ws.Range("H151").Value = "'789500.00"
$ws.Range("H152").Value = "'2975.00"
$ws.Range("H153").Value = "'97800.00"
$ws.Range("H154").Value = "'5432.65"
$ws.Range("H155").Value = "'4560.00"
$ws.Range("H156").Value = "'5328.40"
$ws.Range("H157").Value = "'91154.96"
$ws.Range("H158").Value = "'10770.00"
$ws.Range("H159").Value = "'17275.00"
$ws.Range("H160").Value = "'11545.00"
$ws.Range("H161").Value = "'3881.48"
$ws.Range("H162").Value = "'10000000.00"
$ws.Range("H163").Value = "'60925.00"
$ws.Range("H164").Value = "'16000.00"
$ws.Range("H165").Value = "'30500.00"
$ws.Range("H166").Value = "'22000.00"
$ws.Range("H167").Value = "'40000.00"
$ws.Range("H168").Value = "'14000.00"
$ws.Range("H169").Value = "'22000.00"
$ws.Range("H170").Value = "'18500.00"
$ws.Range("H171").Value = "'15500.00"
$ws.Range("H172").Value = "'22000.25"
$ws.Range("H173").Value = "'21500.00"
$ws.Range("H174").Value = "'18000.00"
$ws.Range("H175").Value = "'16000.00"
$ws.Range("H176").Value = "'16000.00"
$ws.Range("H177").Value = "'18000.00"
$ws.Range("H178").Value = "'18000.00"
$ws.Range("H179").Value = "'10000.00"
$ws.Range("H180").Value = "'7000.00"
$ws.Range("H181").Value = "'31000.00"
$ws.Range("H182").Value = "'18000.00"
$ws.Range("H183").Value = "'18000.00"
$ws.Range("H184").Value = "'18500.00"
$ws.Range("H185").Value = "'15000.00"
$ws.Range("H186").Value = "'16000.00"
$ws.Range("H187").Value = "'18000.00"
$ws.Range("H188").Value = "'22000.00"
$ws.Range("H189").Value = "'8000.00"
$ws.Range("H190").Value = "'22800.00"
$ws.Range("H191").Value = "'2500.00"
$ws.Range("H192").Value = "'50573.21"
$ws.Range("H193").Value = "'33000.00"
$ws.Range("H194").Value = "'18500.00"
$ws.Range("H195").Value = "'18000.00"
$ws.Range("H196").Value = "'41278.55"
$ws.Range("H197").Value = "'7000.00"
$ws.Range("H198").Value = "'90780.00"
$ws.Range("H199").Value = "'18000.00"
$ws.Range("H200").Value = "'12800.00"
$ws.Range("H201").Value = "'4115.23"
$ws.Range("H202").Value = "'12800.00"
$ws.Range("H203").Value = "'650.02"
$ws.Range("H204").Value = "'9565.00"
$ws.Range("H205").Value = "'6655.00"
$ws.Range("H206").Value = "'2200.00"
$ws.Range("H207").Value = "'23710.00"
$ws.Range("H208").Value = "'16345.00"
$ws.Range("H209").Value = "'175.00"
$ws.Range("H210").Value = "'530.00"
$ws.Range("H211").Value = "'3056.82"
$ws.Range("H212").Value = "'5923.48"
$ws.Range("H213").Value = "'38905.20"
$ws.Range("H214").Value = "'11270.00"
$ws.Range("H215").Value = "'13481.12"
$ws.Range("H216").Value = "'2500.00"
$ws.Range("H217").Value = "'42080.00"
$ws.Range("H218").Value = "'26468.00"
$ws.Range("H219").Value = "'17569.66"
$ws.Range("H220").Value = "'3890.00"
$ws.Range("H221").Value = "'7229.73"
$ws.Range("H222").Value = "'9450.00"
$ws.Range("H223").Value = "'299.88"
$ws.Range("H224").Value = "'4767.00"
$ws.Range("H225").Value = "'943.58"
$ws.Range("H226").Value = "'15030.00"
$ws.Range("H227").Value = "'13170.00"
$ws.Range("H228").Value = "'2000.00"
$ws.Range("H229").Value = "'4615.74"
$ws.Range("H230").Value = "'15754.20"
$ws.Range("H231").Value = "'60000.00"
$ws.Range("H232").Value = "'30000.00"
$ws.Range("H233").Value = "'30000.00"
$ws.Range("H234").Value = "'60000.00"
$ws.Range("H235").Value = "'30000.00"
$ws.Range("H236").Value = "'45000.00"
$ws.Range("H237").Value = "'30000.00"
$ws.Range("H238").Value = "'30000.00"
$ws.Range("H239").Value = "'60000.00"
$ws.Range("H240").Value = "'60000.00"
$ws.Range("H241").Value = "'1461.18"
$ws.Range("H242").Value = "'19500.00"
$ws.Range("H243").Value = "'850.00"
$ws.Range("H244").Value = "'17976.95"
$ws.Range("H245").Value = "'4942777.96"
$ws.Range("H246").Value = "'9550.00"
$ws.Range("H247").Value = "'738500.00"
$ws.Range("H248").Value = "'9000.00"
$ws.Range("H249").Value = "'341167.30"
$ws.Range("H250").Value = "'241000.00"
$ws.Range("H251").Value = "'22944976.12"
$ws.Range("H252").Value = "'325800.00"
$ws.Range("H253").Value = "'241000.00"
$ws.Range("H254").Value = "'254200.00"
$ws.Range("H255").Value = "'241000.00"
$ws.Range("H256").Value = "'241000.00"
$ws.Range("H257").Value = "'451000.00"
$ws.Range("H258").Value = "'241000.00"
$ws.Range("H259").Value = "'366900.00"
$ws.Range("H260").Value = "'513000.00"
$ws.Range("H261").Value = "'337650.00"
$ws.Range("H262").Value = "'241000.00"
$ws.Range("H263").Value = "'241000.00"
$ws.Range("H264").Value = "'482000.00"
$ws.Range("H265").Value = "'394150.00"
$ws.Range("H266").Value = "'474400.00"
$ws.Range("H267").Value = "'929050.00"
$ws.Range("H268").Value = "'451000.00"
$ws.Range("H269").Value = "'745900.00"
$ws.Range("H270").Value = "'482000.00"
$ws.Range("H271").Value = "'257410.00"
$ws.Range("H272").Value = "'118459.50"
$ws.Range("H273").Value = "'92202.00"
$ws.Range("H274").Value = "'26018693.51"
$ws.Range("H275").Value = "'1852723.99"
$ws.Range("H276").Value = "'27930.00"
$ws.Range("H277").Value = "'639995.34"
$ws.Range("H278").Value = "'1280000.00"
$ws.Range("H279").Value = "'250000.00"
$ws.Range("H280").Value = "'1500.00"
$ws.Range("H281").Value = "'118000.00"
$ws.Range("H282").Value = "'14500.00"
$ws.Range("H283").Value = "'9600.00"
$ws.Range("H284").Value = "'15212.50"
$ws.Range("H285").Value = "'8000.00"
$ws.Range("H286").Value = "'48000.00"
$ws.Range("H287").Value = "'6090.00"
$ws.Range("H288").Value = "'40900.00"
$ws.Range("H289").Value = "'744000.00"
$ws.Range("H290").Value = "'17000.00"
$ws.Range("H291").Value = "'11055.00"
$ws.Range("H292").Value = "'5500.00"
$ws.Range("H293").Value = "'22350.00"

# Restore original styles (removes any quote-prefix flag picked up above)
$hRange.Style = $hStyle
$ws.Range("E205").Style = $s_E205
$ws.Range("E207").Style = $s_E207
$ws.Range("E217").Style = $s_E217

Write-Host "Updated" 292 "Importe cells and" 3 "Razon social cells"
